$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-21 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-22 Thursday", 2) | Out-Null
$d.Content.Find.Execute("60×63=", $true, $false, $false, $false, $false, $true, 1, $false, "41×54=", 2) | Out-Null
$d.Content.Find.Execute("87×18=", $true, $false, $false, $false, $false, $true, 1, $false, "60×31=", 2) | Out-Null
$d.Content.Find.Execute("35×99=", $true, $false, $false, $false, $false, $true, 1, $false, "29×39=", 2) | Out-Null
$d.Content.Find.Execute("65×50=", $true, $false, $false, $false, $false, $true, 1, $false, "55×36=", 2) | Out-Null
$d.Content.Find.Execute("95×25=", $true, $false, $false, $false, $false, $true, 1, $false, "30×65=", 2) | Out-Null
$d.Content.Find.Execute("58×83=", $true, $false, $false, $false, $false, $true, 1, $false, "34×86=", 2) | Out-Null
$d.Content.Find.Execute("18×50=", $true, $false, $false, $false, $false, $true, 1, $false, "68×79=", 2) | Out-Null
$d.Content.Find.Execute("52×53=", $true, $false, $false, $false, $false, $true, 1, $false, "69×13=", 2) | Out-Null
$d.Content.Find.Execute("15×62=", $true, $false, $false, $false, $false, $true, 1, $false, "59×95=", 2) | Out-Null
$d.Content.Find.Execute("51×84=", $true, $false, $false, $false, $false, $true, 1, $false, "90×96=", 2) | Out-Null
$d.Content.Find.Execute("61×54=", $true, $false, $false, $false, $false, $true, 1, $false, "83×53=", 2) | Out-Null
$d.Content.Find.Execute("24×98=", $true, $false, $false, $false, $false, $true, 1, $false, "13×52=", 2) | Out-Null
$d.Content.Find.Execute("69×81=", $true, $false, $false, $false, $false, $true, 1, $false, "39×74=", 2) | Out-Null
$d.Content.Find.Execute("96×11=", $true, $false, $false, $false, $false, $true, 1, $false, "98×94=", 2) | Out-Null
$d.Content.Find.Execute("84×60=", $true, $false, $false, $false, $false, $true, 1, $false, "59×57=", 2) | Out-Null
$d.Content.Find.Execute("19×51=", $true, $false, $false, $false, $false, $true, 1, $false, "63×92=", 2) | Out-Null
$d.Content.Find.Execute("94×26=", $true, $false, $false, $false, $false, $true, 1, $false, "78×28=", 2) | Out-Null
$d.Content.Find.Execute("24×61=", $true, $false, $false, $false, $false, $true, 1, $false, "64×39=", 2) | Out-Null
$d.Content.Find.Execute("73×54=", $true, $false, $false, $false, $false, $true, 1, $false, "50×66=", 2) | Out-Null
$d.Content.Find.Execute("26×60=", $true, $false, $false, $false, $false, $true, 1, $false, "89×68=", 2) | Out-Null
$d.Content.Find.Execute("96×95=", $true, $false, $false, $false, $false, $true, 1, $false, "94×87=", 2) | Out-Null
$d.Content.Find.Execute("29×68=", $true, $false, $false, $false, $false, $true, 1, $false, "57×93=", 2) | Out-Null
$d.Content.Find.Execute("15×31=", $true, $false, $false, $false, $false, $true, 1, $false, "79×90=", 2) | Out-Null
$d.Content.Find.Execute("89×55=", $true, $false, $false, $false, $false, $true, 1, $false, "35×57=", 2) | Out-Null
